$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9527859091758728
$ws.Range("B1").Value = 1.696652889251709
$ws.Range("C1").Value = 5.610214710235596
$ws.Range("D1").Value = 3.572016477584839
$ws.Range("E1").Value = 1.120502114295959
